$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: lname Walsh -> Patterson2
$ws.Range("C4").Value = "Patterson2"

# Row 5: lname Patterson -> Pattersons, contact Contact@test.com -> Contact@test.cool
$ws.Range("C5").Value = "Pattersons"
$ws.Range("E5").Value = "Contact@test.cool"

# Row 6: fname Akiel -> Jessa, lname Jame -> James, company Atlas -> Atlast, contact jessaj@test.com -> jj@test.works
$ws.Range("B6").Value = "Jessa"
$ws.Range("C6").Value = "James"
$ws.Range("D6").Value = "Atlast"
$ws.Range("E6").Value = "jj@test.works"
